$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 1588
$ws.Range("I2").Value = 4550
$ws.Range("J2").Value = 18640
$ws.Range("K2").Value = 87
$ws.Range("L2").Value = 5145
$ws.Range("M2").Value = 326
$ws.Range("N2").Value = 3108
$ws.Range("P2").Value = 76
$ws.Range("Q2").Value = 28
$ws.Range("R2").Value = 282
$ws.Range("S2").Value = 2021
$ws.Range("T2").Value = 3241
$ws.Range("U2").Value = 248
$ws.Range("V2").Value = 28957
$ws.Range("W2").Value = 12
$ws.Range("X2").Value = 28840
$ws.Range("Y2").Value = 44
$ws.Range("Z2").Value = 413
$ws.Range("AA2").Value = 206
